$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $n = $parts.Length
            $revParts = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $revParts += $parts[$i]
            }
            $newVal = $revParts -join ", "
            $cell.Value = $newVal
        }
    }
}
